# Insert a new weekly record row before the current row 280, shifting the
# existing rows 280-389 down to 281-390 (row 389's old content lands in the
# new row 390). Then populate the newly inserted row 280 with the new
# record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 280 (pushes rows 280..389 down to 281..390).
$ws.Rows.Item(280).Insert()

# Fill in the new row 280 with the new weekly record.
$ws.Cells.Item(280, 1).Value = 9
$ws.Cells.Item(280, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(280, 3).Value = 'Metropolitana'
$ws.Cells.Item(280, 4).Value = 44726
$ws.Cells.Item(280, 5).Value = 13
$ws.Cells.Item(280, 6).Value = 100112044
$ws.Cells.Item(280, 7).Value = 'Perejil'
$ws.Cells.Item(280, 8).Value = 'Sin especificar'
$ws.Cells.Item(280, 9).Value = 'Primera'
$ws.Cells.Item(280, 10).Value = 93
$ws.Cells.Item(280, 11).Value = 8000
$ws.Cells.Item(280, 12).Value = 10000
$ws.Cells.Item(280, 13).Value = 9075
$ws.Cells.Item(280, 14).Value = '$/docena de atados'
$ws.Cells.Item(280, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(280, 16).Value = 3025
$ws.Cells.Item(280, 17).Value = 3
$ws.Cells.Item(280, 18).Value = 'Hortaliza'
